# Applies the "python run_script and associated files" commit to
# misc/logos_for_README.pptx:
#   - slide 4, group "Group 20" / shape "TextBox 4":
#       * reposition / resize the textbox
#       * append ".fq.gz" to the sample-name text
#
# PowerPoint's COM object model works in points, while the OOXML stores
# EMUs (1 pt = 12700 EMU). EMU-per-point helper below also nudges the
# value a half-EMU "up" before the pt<->EMU round trip; the host stores
# Left/Top/Width/Height internally as 32-bit floats and truncates (not
# rounds) when it converts back to EMU on save, so without the nudge a
# plain emu/12700.0 can land one EMU short after that truncation.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$grp = $s.Shapes.Item(1)

$shp = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($item.Name -eq "TextBox 4") {
        $shp = $item
    }
}

# Update the text first -- this shape auto-fits to its text
# (<a:spAutoFit/>), so the explicit size below must be applied after
# the text change or the auto-fit recalculation clobbers it.
$shp.TextFrame.TextRange.Text = "S101-F010-L01-S101F010L01.fq.gz"

$shp.Left = EmuToPt 2522400
$shp.Top = EmuToPt 2947855
$shp.Width = EmuToPt 4458272
$shp.Height = EmuToPt 369332
